$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.880.31"
$ws.Range('E2').Value = '  -0.54%  '

$ws.Range('D3').Value = "'1.808.17"
$ws.Range('E3').Value = '  -1.21%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = "'240.36"
$ws.Range('E5').Value = '  -1.28%  '

$ws.Range('D6').Value = "'0.6053"
$ws.Range('E6').Value = '  -3.62%  '

$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = "'0.07254"
$ws.Range('E8').Value = '  -3.11%  '

$ws.Range('D9').Value = "'0.2856"
$ws.Range('E9').Value = '  -2.15%  '

$ws.Range('D10').Value = "'22.68"
$ws.Range('E10').Value = '  -2.01%  '

$ws.Range('D11').Value = "'0.07620"
$ws.Range('E11').Value = '  -1.19%  '

$ws.Range('D12').Value = "'1.821.60"
$ws.Range('E12').Value = '  -0.56%  '

$ws.Range('D13').Value = "'4.905"
$ws.Range('E13').Value = '  -1.67%  '

$ws.Range('D14').Value = "'0.6542"
$ws.Range('E14').Value = '  -2.00%  '

$ws.Range('D15').Value = "'80.69"
$ws.Range('E15').Value = '  -2.26%  '

$ws.Range('D16').Value = "'0.000008900"
$ws.Range('E16').Value = '  -4.58%  '

$ws.Range('D17').Value = "'5.809"
$ws.Range('E17').Value = '  -3.06%  '

$ws.Range('D18').Value = "'28.870.05"
$ws.Range('E18').Value = '  -0.67%  '

$ws.Range('D19').Value = "'2.061.96"
$ws.Range('E19').Value = '  -0.64%  '

$ws.Range('D20').Value = "'235.72"
$ws.Range('E20').Value = '  +5.51%  '

$ws.Range('D21').Value = "'12.33"
$ws.Range('E21').Value = '  -2.01%  '

$ws.Range('D22').Value = "'1.000"
$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('D23').Value = "'7.062"
$ws.Range('E23').Value = '  -0.70%  '

$ws.Range('D24').Value = "'1.001"
$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').Value = "'157.88"
$ws.Range('E25').Value = '  -0.85%  '

$ws.Range('D26').Value = "'0.1392"
$ws.Range('E26').Value = '  -0.39%  '

$ws.Range('D27').Value = "'8.346"
$ws.Range('E27').Value = '  -1.66%  '

$ws.Range('D28').Value = "'17.49"
$ws.Range('E28').Value = '  -2.40%  '

$ws.Range('D29').Value = "'1.473"
$ws.Range('E29').Value = '  -1.48%  '

$ws.Range('D30').Value = "'0.05563"
$ws.Range('E30').Value = '  -2.38%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'4.056"
$ws.Range('E31').Value = '  -2.28%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'4.038"
$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('D33').Value = "'1.204"
$ws.Range('E33').Value = '  +0.20%  '

$ws.Range('D34').Value = "'1.804"
$ws.Range('E34').Value = '  -2.10%  '

$ws.Range('D35').Value = "'0.7280"
$ws.Range('E35').Value = '  -2.51%  '

$ws.Range('D36').Value = "'1.122"
$ws.Range('E36').Value = '  -1.09%  '

$ws.Range('D37').Value = "'2.626"
$ws.Range('E37').Value = '  -1.48%  '

$ws.Range('D38').Value = "'2.810"
$ws.Range('E38').Value = '  +1.92%  '

$ws.Range('D39').Value = "'0.01740"
$ws.Range('E39').Value = '  -2.39%  '

$ws.Range('D40').Value = "'1.185.30"
$ws.Range('E40').Value = '  -2.56%  '

$ws.Range('D41').Value = "'6.327"
$ws.Range('E41').Value = '  -3.25%  '

$ws.Range('D42').Value = "'0.8852"
$ws.Range('E42').Value = '  -0.96%  '

$ws.Range('D43').Value = "'0.9999"
$ws.Range('E43').Value = '  +0.04%  '

$ws.Range('D44').Value = "'99.98"
$ws.Range('E44').Value = '  -2.07%  '

$ws.Range('D45').Value = "'1.964.90"
$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000123"
$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'63.84"
$ws.Range('E47').Value = '  -2.60%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.5076"
$ws.Range('E48').Value = '  -0.03%  '

$ws.Range('D49').Value = "'8.994"
$ws.Range('E49').Value = '  -0.25%  '

$ws.Range('D50').Value = "'0.3952"
$ws.Range('E50').Value = '  -2.85%  '

$ws.Range('D51').Value = "'0.05769"
$ws.Range('E51').Value = '  -0.91%  '
